$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-24 22:48:18"
$ws.Range("O2").Value = "5.6 °C"
$ws.Range("E3").Value = "2026-02-24 22:48:21"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "33%"
$ws.Range("E4").Value = "2026-02-24 22:48:23"
$ws.Range("O4").Value = "12.4 °C"
$ws.Range("E5").Value = "2026-02-24 22:48:26"
$ws.Range("E6").Value = "2026-02-24 22:48:28"
$ws.Range("E7").Value = "2026-02-24 22:48:30"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "74%"
$ws.Range("J7").Value = "1020.0 hPa"
$ws.Range("E8").Value = "2026-02-24 22:48:33"
$ws.Range("E9").Value = "2026-02-24 22:48:35"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "83%"
$ws.Range("O9").Value = "11.6 °C"
$ws.Range("E10").Value = "2026-02-24 22:48:38"
$ws.Range("O10").Value = "10.8 °C"
$ws.Range("E11").Value = "2026-02-24 22:48:40"
$ws.Range("O11").Value = "8.9 °C"
$ws.Range("E12").Value = "2026-02-24 22:48:43"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "92%"
$ws.Range("O12").Value = "10.3 °C"
$ws.Range("E13").Value = "2026-02-24 22:48:44"
$ws.Range("J13").Value = "1023.3 hPa"
$ws.Range("O13").Value = "6.6 °C"
$ws.Range("E14").Value = "2026-02-24 22:48:47"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "86%"
$ws.Range("N14").Value = "5.3 °C 22:04 TU"
$ws.Range("O14").Value = "11.2 °C"
$ws.Range("E15").Value = "2026-02-24 22:48:49"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "78%"
$ws.Range("O15").Value = "11.6 °C"
$ws.Range("E16").Value = "2026-02-24 22:48:52"
$ws.Range("E17").Value = "2026-02-24 22:48:54"
$ws.Range("E18").Value = "2026-02-24 22:48:57"
$ws.Range("O18").Value = "11.0 °C"
$ws.Range("E19").Value = "2026-02-24 22:48:59"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "52%"
$ws.Range("E20").Value = "2026-02-24 22:49:01"
$ws.Range("E21").Value = "2026-02-24 22:49:04"
$ws.Range("J21").Value = "1021.8 hPa"
$ws.Range("E22").Value = "2026-02-24 22:49:06"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "27%"
$ws.Range("L22").Value = "23.0 km/h - 316º 22:01 TU"
$ws.Range("N22").Value = "0.9 °C 22:18 TU"
$ws.Range("E23").Value = "2026-02-24 22:49:09"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "24%"
$ws.Range("N23").Value = "2.8 °C 22:01 TU"
$ws.Range("E24").Value = "2026-02-24 22:49:11"
$ws.Range("E25").Value = "2026-02-24 22:49:14"
$ws.Range("N25").Value = "3.1 °C 22:24 TU"
$ws.Range("O25").Value = "6.5 °C"
$ws.Range("E26").Value = "2026-02-24 22:49:16"
$ws.Range("J26").Value = "1018.7 hPa"
$ws.Range("E27").Value = "2026-02-24 22:49:19"
$ws.Range("L27").Value = "17.6 km/h - 202º 22:07 TU"
$ws.Range("O27").Value = "6.2 °C"
$ws.Range("E28").Value = "2026-02-24 22:49:21"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "72%"
$ws.Range("O28").Value = "11.3 °C"
$ws.Range("E29").Value = "2026-02-24 22:49:24"
$ws.Range("O29").Value = "10.0 °C"
$ws.Range("E30").Value = "2026-02-24 22:49:26"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "77%"
$ws.Range("O30").Value = "12.9 °C"
$ws.Range("E31").Value = "2026-02-24 22:49:29"
$ws.Range("K31").Value = "15.1 MJ/m2"
$ws.Range("N31").Value = "12.3 °C 22:15 TU"
$ws.Range("E32").Value = "2026-02-24 22:49:31"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "70%"
$ws.Range("O32").Value = "6.8 °C"
$ws.Range("E33").Value = "2026-02-24 22:49:33"
$ws.Range("E34").Value = "2026-02-24 22:49:36"
$ws.Range("O34").Value = "4.5 °C"
$ws.Range("E35").Value = "2026-02-24 22:49:38"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = "35%"
$ws.Range("E36").Value = "2026-02-24 22:49:41"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "82%"
$ws.Range("J36").Value = "1019.7 hPa"
$ws.Range("E37").Value = "2026-02-24 22:49:43"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "73%"
$ws.Range("J37").Value = "1022.4 hPa"
$ws.Range("O37").Value = "8.4 °C"
$ws.Range("E38").Value = "2026-02-24 22:49:45"
$ws.Range("O38").Value = "11.7 °C"
$ws.Range("E39").Value = "2026-02-24 22:49:48"
$ws.Range("N39").Value = "1.0 °C 22:25 TU"
$ws.Range("O39").Value = "4.4 °C"
$ws.Range("E40").Value = "2026-02-24 22:49:50"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "67%"
$ws.Range("J40").Value = "1022.7 hPa"
$ws.Range("O40").Value = "8.3 °C"
$ws.Range("E41").Value = "2026-02-24 22:49:53"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "81%"
$ws.Range("E42").Value = "2026-02-24 22:49:55"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "88%"
$ws.Range("O42").Value = "11.1 °C"
$ws.Range("E43").Value = "2026-02-24 22:49:57"
$ws.Range("O43").Value = "10.6 °C"
$ws.Range("E44").Value = "2026-02-24 22:50:00"
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = "41%"
$ws.Range("E45").Value = "2026-02-24 22:50:02"
$ws.Range("K45").Value = "13.8 MJ/m2"
$ws.Range("O45").Value = "9.9 °C"
$ws.Range("E46").Value = "2026-02-24 22:50:05"
$ws.Range("O46").Value = "10.6 °C"
